$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the formulas across each column range so Excel groups them as shared formulas
$ws.Range("I2:I14").Formula = "=ROUND(IMDIV(B2,H2),2)"
$ws.Range("J2:J14").Formula = "=ROUND(IMDIV(D2,H2),2)"
$ws.Range("K2:K14").Formula = "=ROUND(IMDIV(E2,H2),2)"
$ws.Range("L2:L14").Formula = "=ROUND(IMDIV(F2,H2),2)"
$ws.Range("M2:M14").Formula = "=ROUND(IMDIV(B2,D2),2)"
$ws.Range("N2:N14").Formula = "=ROUND(IMDIV(C2,B2),4)*100"
$ws.Range("O2:O14").Formula = "=ROUND(IMDIV(G2,H2),2)"

# Fill in "none" for the Awards column where it was previously blank
$ws.Range("Q3").Value = "none"
$ws.Range("Q5").Value = "none"
$ws.Range("Q6").Value = "none"
$ws.Range("Q9").Value = "none"
$ws.Range("Q10").Value = "none"
$ws.Range("Q11").Value = "none"
$ws.Range("Q12").Value = "none"
$ws.Range("Q13").Value = "none"
$ws.Range("Q14").Value = "none"

# Update the active selection/view to reflect the saved cursor position
$ws.Range("U20").Select() | Out-Null
